# Remove the second table (the small "Men"/"Women" summary table that
# follows the main demographics table). The rest of the document is
# left untouched.
$d = $word.ActiveDocument

if ($d.Tables.Count -ge 2) {
    $t = $d.Tables.Item(2)
    $t.Delete()
}
